# Generate Report for Handoff
# - Flip the localization status from "In Translation" to "Ready for handoff"
#   on the Overview sheet (per-language status columns) and on each
#   per-language detail sheet ("zh-cn", "de-de").
# - Refresh the associated timestamps (Latest HO Xliff Generate Date /
#   Latest Handoff Datetime) to the new handoff moment.
# - The Status columns now hold a longer string, so widen them to fit
#   (matches the column-width bump baked into the target workbook).

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus            # zh-cn status
$wsOverview.Range("F2").Value = $newStatus            # de-de status
$wsOverview.Range("G2").Value = "2016-08-31 04:40:58" # Latest HO Xliff Generate Date

$wsOverview.Columns.Item(5).ColumnWidth = 16.38
$wsOverview.Columns.Item(6).ColumnWidth = 16.38

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus                # Status
$wsZhCn.Range("H2").Value = "2016-08-31 04:40:53"     # Latest Handoff Datetime

$wsZhCn.Columns.Item(3).ColumnWidth = 16.38

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus                # Status
$wsDeDe.Range("H2").Value = "2016-08-31 04:40:58"     # Latest Handoff Datetime

$wsDeDe.Columns.Item(3).ColumnWidth = 16.38
